# Update the 2024-02-06 [N]/Δ[N]/[%]/Δ[%] values (columns X:AA) for the
# "latitude" (row 5) and "longitude" (row 6) fields to reflect the newly
# mapped 1km/10km reference grid counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 5, 6) {
    $ws.Range("X$r").Value = 10735
    $ws.Range("Y$r").Value = -98
    $ws.Range("Z$r").Value = 98.7217215376127
    $ws.Range("AA$r").Value = -0.901232297222734
}
